$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ContactUs")

# Fill in the missing "soft" value for the Developer row (row 3), matching
# the pattern of Softway23/Softway26/Softway25 in column C of the other rows.
$ws.Range("C3").Value = "soft"
$ws.Range("C3").WrapText = $true
